$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the table by two rows: duplicate formatting of the last existing
#     row (15) down into the two brand-new rows (16, 17) so the new cells
#     pick up the same style (border/bold) as the rest of column A and the
#     same value-types (string / boolean) as columns B and E.
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A15:E15").Copy($ws.Range("A17:E17"))

# --- Row 8: was extr1, now line7 ---------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $True

# --- Row 9: was extr2, now line8 ---------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $True

# --- Row 10: was extr3, now extr1 ---------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $True

# --- Row 11: was extr4, now extr2 ---------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $True

# --- Row 12: was extr5, now extr3 ---------------------------------------
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $True

# --- Row 13: was extr6, now extr4 ---------------------------------------
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $True

# --- Row 14: was extr7, now extr5 ---------------------------------------
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $True

# --- Row 15: was extr8, now extr6 ---------------------------------------
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $True

# --- Row 16 (new): extr7 -------------------------------------------------
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $True

# --- Row 17 (new): extr8 -------------------------------------------------
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $True
